{"js": "// Apply the text replacements described by the diff: the header date and\n// each \"NNN\u00f7N=\" exercise cell get replaced with new values, in document\n// order. All \"old\" strings are unique in the document, so a straightforward\n// search + replace for each pair is safe.\nconst replacements = [\n  [\"2025-08-29 Friday\", \"2025-08-30 Saturday\"],\n  [\"419\u00f72=\", \"999\u00f74=\"],\n  [\"574\u00f75=\", \"829\u00f75=\"],\n  [\"783\u00f79=\", \"793\u00f79=\"],\n  [\"878\u00f76=\", \"620\u00f77=\"],\n  [\"637\u00f72=\", \"249\u00f76=\"],\n  [\"677\u00f74=\", \"132\u00f75=\"],\n  [\"613\u00f73=\", \"623\u00f78=\"],\n  [\"500\u00f76=\", \"528\u00f79=\"],\n  [\"167\u00f79=\", \"465\u00f75=\"],\n  [\"337\u00f79=\", \"335\u00f78=\"],\n  [\"359\u00f73=\", \"313\u00f76=\"],\n  [\"503\u00f75=\", \"956\u00f77=\"],\n  [\"442\u00f77=\", \"442\u00f74=\"],\n  [\"252\u00f77=\", \"227\u00f79=\"],\n  [\"784\u00f73=\", \"471\u00f76=\"],\n  [\"192\u00f76=\", \"267\u00f79=\"],\n  [\"117\u00f79=\", \"858\u00f79=\"],\n  [\"702\u00f79=\", \"583\u00f73=\"],\n  [\"257\u00f72=\", \"564\u00f72=\"],\n  [\"856\u00f72=\", \"824\u00f76=\"],\n  [\"654\u00f77=\", \"390\u00f76=\"],\n  [\"266\u00f78=\", \"290\u00f79=\"],\n  [\"480\u00f73=\", \"580\u00f72=\"],\n  [\"320\u00f76=\", \"556\u00f73=\"],\n  [\"257\u00f78=\", \"303\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the text replacements described by the diff: the header date and\n# each \"NNN\u00f7N=\" exercise cell get replaced with new values. All \"old\"\n# strings are unique in the document, so a plain Find/Replace for each\n# pair (wildcards off) is safe and order-independent.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-08-29 Friday\", \"2025-08-30 Saturday\"),\n    @(\"419\u00f72=\", \"999\u00f74=\"),\n    @(\"574\u00f75=\", \"829\u00f75=\"),\n    @(\"783\u00f79=\", \"793\u00f79=\"),\n    @(\"878\u00f76=\", \"620\u00f77=\"),\n    @(\"637\u00f72=\", \"249\u00f76=\"),\n    @(\"677\u00f74=\", \"132\u00f75=\"),\n    @(\"613\u00f73=\", \"623\u00f78=\"),\n    @(\"500\u00f76=\", \"528\u00f79=\"),\n    @(\"167\u00f79=\", \"465\u00f75=\"),\n    @(\"337\u00f79=\", \"335\u00f78=\"),\n    @(\"359\u00f73=\", \"313\u00f76=\"),\n    @(\"503\u00f75=\", \"956\u00f77=\"),\n    @(\"442\u00f77=\", \"442\u00f74=\"),\n    @(\"252\u00f77=\", \"227\u00f79=\"),\n    @(\"784\u00f73=\", \"471\u00f76=\"),\n    @(\"192\u00f76=\", \"267\u00f79=\"),\n    @(\"117\u00f79=\", \"858\u00f79=\"),\n    @(\"702\u00f79=\", \"583\u00f73=\"),\n    @(\"257\u00f72=\", \"564\u00f72=\"),\n    @(\"856\u00f72=\", \"824\u00f76=\"),\n    @(\"654\u00f77=\", \"390\u00f76=\"),\n    @(\"266\u00f78=\", \"290\u00f79=\"),\n    @(\"480\u00f73=\", \"580\u00f72=\"),\n    @(\"320\u00f76=\", \"556\u00f73=\"),\n    @(\"257\u00f78=\", \"303\u00f74=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
